$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 337, shifting existing rows 337:358 down to 338:359.
$ws.Rows(337).Insert()

# Populate the newly inserted row 337 with the new record.
$ws.Range("A337").Value = 4
$ws.Range("B337").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C337").Value = "Los Lagos"
$ws.Range("D337").Value = 44610
$ws.Range("E337").Value = 10
$ws.Range("F337").Value = "Fruta"
$ws.Range("G337").Value = 100102
$ws.Range("H337").Value = "Cítricos"
$ws.Range("I337").Value = 100102005
$ws.Range("J337").Value = "Naranja"
$ws.Range("K337").Value = "Valencia"
$ws.Range("L337").Value = "Primera"
$ws.Range("M337").Value = 400
$ws.Range("N337").Value = 16000
$ws.Range("O337").Value = 16000
$ws.Range("P337").Value = 16000
$ws.Range("Q337").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R337").Value = "Región de O'Higgins"
$ws.Range("S337").Value = 1067
$ws.Range("T337").Value = 15
